$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2115
$ws.Range("K3").Value = 2054
$ws.Range("K4").Value = 431
$ws.Range("K5").Value = 137
$ws.Range("K6").Value = 2603
$ws.Range("K7").Value = 7340

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 142
$ws.Range("K3").Value = 141
$ws.Range("K7").Value = 495

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 63
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 82
$ws.Range("K4").Value = 17
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 290

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 80
$ws.Range("K7").Value = 238

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 51
$ws.Range("K3").Value = 43
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 54
$ws.Range("K5").Value = 12
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 215
$ws.Range("K8").Value = 495
$ws.Range("K11").Value = 158
$ws.Range("K15").Value = 71
$ws.Range("K19").Value = 207
$ws.Range("K20").Value = 155
$ws.Range("K29").Value = 365
$ws.Range("K33").Value = 290
$ws.Range("K35").Value = 12
$ws.Range("K36").Value = 84
$ws.Range("K37").Value = 238
$ws.Range("K42").Value = 255
$ws.Range("K43").Value = 67
$ws.Range("K46").Value = 15
$ws.Range("K47").Value = 43
$ws.Range("K48").Value = 91
$ws.Range("K49").Value = 50
$ws.Range("K52").Value = 196
$ws.Range("K54").Value = 132
$ws.Range("K55").Value = 78
$ws.Range("K63").Value = 26
$ws.Range("K64").Value = 48
$ws.Range("K65").Value = 178
$ws.Range("K67").Value = 280
$ws.Range("K73").Value = 73
$ws.Range("K77").Value = 52
$ws.Range("K78").Value = 96
$ws.Range("K79").Value = 196
$ws.Range("K83").Value = 156
$ws.Range("K85").Value = 360
$ws.Range("K86").Value = 49
$ws.Range("K90").Value = 62
$ws.Range("K91").Value = 69
$ws.Range("K95").Value = 114
$ws.Range("K96").Value = 102
$ws.Range("K99").Value = 135
$ws.Range("K101").Value = 7340

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 80
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 94
$ws.Range("K4").Value = 20
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 365

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 207

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 60
$ws.Range("K3").Value = 76
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 29
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 28
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 53
$ws.Range("K7").Value = 155

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 73
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 17
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 158

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K4").Value = 1
$ws.Range("K6").Value = 12

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 10
$ws.Range("K6").Value = 49

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 130
$ws.Range("K3").Value = 121
$ws.Range("K7").Value = 360

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 43
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 196
